# Auto-generated: rebuild Sheet1 data rows 2-21 (Tgfb2-Tgfbr2 NATMI LR-pairs TPM update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 20,20

$data[0,0] = "ECs"
$data[0,1] = "Tgfb2"
$data[0,2] = "Tgfbr2"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 2.767474
$data[0,7] = 8.302422
$data[0,8] = 0.1192484523516842
$data[0,9] = 0.1192484523516842
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 42.99144133333333
$data[0,13] = 128.974324
$data[0,14] = 0.2509605789056467
$data[0,15] = 0.2509605789056467
$data[0,16] = 118.9776961125253
$data[0,17] = 1070.799265012728
$data[0,18] = 0.02992666063578108
$data[0,19] = 0.02992666063578109

$data[1,0] = "ECs"
$data[1,1] = "Tgfb2"
$data[1,2] = "Tgfbr2"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 2.767474
$data[1,7] = 8.302422
$data[1,8] = 0.1192484523516842
$data[1,9] = 0.1192484523516842
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 55.607043
$data[1,13] = 166.821129
$data[1,14] = 0.3246035785195009
$data[1,15] = 0.324603578519501
$data[1,16] = 153.891045719382
$data[1,17] = 1385.019411474438
$data[1,18] = 0.03870847436626888
$data[1,19] = 0.03870847436626889

$data[2,0] = "ECs"
$data[2,1] = "Tgfb2"
$data[2,2] = "Tgfbr2"
$data[2,3] = "Inflammatory-Mac"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 2.767474
$data[2,7] = 8.302422
$data[2,8] = 0.1192484523516842
$data[2,9] = 0.1192484523516842
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 36.72715
$data[2,13] = 110.18145
$data[2,14] = 0.2143930638214748
$data[2,15] = 0.2143930638214748
$data[2,16] = 101.6414327191
$data[2,17] = 914.7728944719001
$data[2,18] = 0.02556604105564671
$data[2,19] = 0.02556604105564672

$data[3,0] = "ECs"
$data[3,1] = "Tgfb2"
$data[3,2] = "Tgfbr2"
$data[3,3] = "MuSCs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 2.767474
$data[3,7] = 8.302422
$data[3,8] = 0.1192484523516842
$data[3,9] = 0.1192484523516842
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 7.381512666666667
$data[3,13] = 22.144538
$data[3,14] = 0.04308924368603855
$data[3,15] = 0.04308924368603857
$data[3,16] = 20.42814438567067
$data[3,17] = 183.853299471036
$data[3,18] = 0.005138325622564676
$data[3,19] = 0.005138325622564678

$data[4,0] = "ECs"
$data[4,1] = "Tgfb2"
$data[4,2] = "Tgfbr2"
$data[4,3] = "Resolving-Mac"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 2.767474
$data[4,7] = 8.302422
$data[4,8] = 0.1192484523516842
$data[4,9] = 0.1192484523516842
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 28.60040066666667
$data[4,13] = 85.801202
$data[4,14] = 0.1669535350673389
$data[4,15] = 0.1669535350673389
$data[4,16] = 79.15086523458267
$data[4,17] = 712.3577871112441
$data[4,18] = 0.01990895067142279
$data[4,19] = 0.0199089506714228

$data[5,0] = "FAPs"
$data[5,1] = "Tgfb2"
$data[5,2] = "Tgfbr2"
$data[5,3] = "ECs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 12.88577866666667
$data[5,7] = 38.657336
$data[5,8] = 0.5552388796954726
$data[5,9] = 0.5552388796954726
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 42.99144133333333
$data[5,13] = 128.974324
$data[5,14] = 0.2509605789056467
$data[5,15] = 0.2509605789056467
$data[5,16] = 553.9781975823182
$data[5,17] = 4985.803778240864
$data[5,18] = 0.1393430706792985
$data[5,19] = 0.1393430706792985

$data[6,0] = "FAPs"
$data[6,1] = "Tgfb2"
$data[6,2] = "Tgfbr2"
$data[6,3] = "FAPs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 12.88577866666667
$data[6,7] = 38.657336
$data[6,8] = 0.5552388796954726
$data[6,9] = 0.5552388796954726
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 55.607043
$data[6,13] = 166.821129
$data[6,14] = 0.3246035785195009
$data[6,15] = 0.324603578519501
$data[6,16] = 716.540048405816
$data[6,17] = 6448.860435652344
$data[6,18] = 0.1802325272823091
$data[6,19] = 0.1802325272823091

$data[7,0] = "FAPs"
$data[7,1] = "Tgfb2"
$data[7,2] = "Tgfbr2"
$data[7,3] = "Inflammatory-Mac"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 12.88577866666667
$data[7,7] = 38.657336
$data[7,8] = 0.5552388796954726
$data[7,9] = 0.5552388796954726
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 36.72715
$data[7,13] = 110.18145
$data[7,14] = 0.2143930638214748
$data[7,15] = 0.2143930638214748
$data[7,16] = 473.2579259574667
$data[7,17] = 4259.3213336172
$data[7,18] = 0.1190393645707156
$data[7,19] = 0.1190393645707156

$data[8,0] = "FAPs"
$data[8,1] = "Tgfb2"
$data[8,2] = "Tgfbr2"
$data[8,3] = "MuSCs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 12.88577866666667
$data[8,7] = 38.657336
$data[8,8] = 0.5552388796954726
$data[8,9] = 0.5552388796954726
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 7.381512666666667
$data[8,13] = 22.144538
$data[8,14] = 0.04308924368603855
$data[8,15] = 0.04308924368603857
$data[8,16] = 95.11653844786311
$data[8,17] = 856.048846030768
$data[8,18] = 0.02392482339116126
$data[8,19] = 0.02392482339116127

$data[9,0] = "FAPs"
$data[9,1] = "Tgfb2"
$data[9,2] = "Tgfbr2"
$data[9,3] = "Resolving-Mac"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 12.88577866666667
$data[9,7] = 38.657336
$data[9,8] = 0.5552388796954726
$data[9,9] = 0.5552388796954726
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 28.60040066666667
$data[9,13] = 85.801202
$data[9,14] = 0.1669535350673389
$data[9,15] = 0.1669535350673389
$data[9,16] = 368.5384327686525
$data[9,17] = 3316.845894917872
$data[9,18] = 0.09269909377198805
$data[9,19] = 0.09269909377198807

$data[10,0] = "Inflammatory-Mac"
$data[10,1] = "Tgfb2"
$data[10,2] = "Tgfbr2"
$data[10,3] = "ECs"
$data[10,4] = 1
$data[10,5] = 0.3333333333333333
$data[10,6] = 0.03423166666666667
$data[10,7] = 0.102695
$data[10,8] = 0.001475017749550216
$data[10,9] = 0.001475017749550216
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 42.99144133333333
$data[10,13] = 128.974324
$data[10,14] = 0.2509605789056467
$data[10,15] = 0.2509605789056467
$data[10,16] = 1.471668689242222
$data[10,17] = 13.24501820318
$data[10,18] = 0.0003701713083232264
$data[10,19] = 0.0003701713083232265

$data[11,0] = "Inflammatory-Mac"
$data[11,1] = "Tgfb2"
$data[11,2] = "Tgfbr2"
$data[11,3] = "FAPs"
$data[11,4] = 1
$data[11,5] = 0.3333333333333333
$data[11,6] = 0.03423166666666667
$data[11,7] = 0.102695
$data[11,8] = 0.001475017749550216
$data[11,9] = 0.001475017749550216
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 55.607043
$data[11,13] = 166.821129
$data[11,14] = 0.3246035785195009
$data[11,15] = 0.324603578519501
$data[11,16] = 1.903521760295
$data[11,17] = 17.131695842655
$data[11,18] = 0.0004787960398837811
$data[11,19] = 0.0004787960398837813

$data[12,0] = "Inflammatory-Mac"
$data[12,1] = "Tgfb2"
$data[12,2] = "Tgfbr2"
$data[12,3] = "Inflammatory-Mac"
$data[12,4] = 1
$data[12,5] = 0.3333333333333333
$data[12,6] = 0.03423166666666667
$data[12,7] = 0.102695
$data[12,8] = 0.001475017749550216
$data[12,9] = 0.001475017749550216
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 36.72715
$data[12,13] = 110.18145
$data[12,14] = 0.2143930638214748
$data[12,15] = 0.2143930638214748
$data[12,16] = 1.257231556416667
$data[12,17] = 11.31508400775
$data[12,18] = 0.0003162335745171276
$data[12,19] = 0.0003162335745171276

$data[13,0] = "Inflammatory-Mac"
$data[13,1] = "Tgfb2"
$data[13,2] = "Tgfbr2"
$data[13,3] = "MuSCs"
$data[13,4] = 1
$data[13,5] = 0.3333333333333333
$data[13,6] = 0.03423166666666667
$data[13,7] = 0.102695
$data[13,8] = 0.001475017749550216
$data[13,9] = 0.001475017749550216
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 7.381512666666667
$data[13,13] = 22.144538
$data[13,14] = 0.04308924368603855
$data[13,15] = 0.04308924368603857
$data[13,16] = 0.2526814811011111
$data[13,17] = 2.27413332991
$data[13,18] = 0.00006355739925160144
$data[13,19] = 0.00006355739925160147

$data[14,0] = "Inflammatory-Mac"
$data[14,1] = "Tgfb2"
$data[14,2] = "Tgfbr2"
$data[14,3] = "Resolving-Mac"
$data[14,4] = 1
$data[14,5] = 0.3333333333333333
$data[14,6] = 0.03423166666666667
$data[14,7] = 0.102695
$data[14,8] = 0.001475017749550216
$data[14,9] = 0.001475017749550216
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 28.60040066666667
$data[14,13] = 85.801202
$data[14,14] = 0.1669535350673389
$data[14,15] = 0.1669535350673389
$data[14,16] = 0.9790393821544445
$data[14,17] = 8.81135443939
$data[14,18] = 0.0002462594275744793
$data[14,19] = 0.0002462594275744794

$data[15,0] = "MuSCs"
$data[15,1] = "Tgfb2"
$data[15,2] = "Tgfbr2"
$data[15,3] = "ECs"
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 7.520146
$data[15,7] = 22.560438
$data[15,8] = 0.324037650203293
$data[15,9] = 0.3240376502032931
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 42.99144133333333
$data[15,13] = 128.974324
$data[15,14] = 0.2509605789056467
$data[15,15] = 0.2509605789056467
$data[15,16] = 323.3019155771013
$data[15,17] = 2909.717240193911
$data[15,18] = 0.08132067628224385
$data[15,19] = 0.08132067628224388

$data[16,0] = "MuSCs"
$data[16,1] = "Tgfb2"
$data[16,2] = "Tgfbr2"
$data[16,3] = "FAPs"
$data[16,4] = 3
$data[16,5] = 1
$data[16,6] = 7.520146
$data[16,7] = 22.560438
$data[16,8] = 0.324037650203293
$data[16,9] = 0.3240376502032931
$data[16,10] = 3
$data[16,11] = 1
$data[16,12] = 55.607043
$data[16,13] = 166.821129
$data[16,14] = 0.3246035785195009
$data[16,15] = 0.324603578519501
$data[16,16] = 418.1730819882779
$data[16,17] = 3763.557737894501
$data[16,18] = 0.1051837808310392
$data[16,19] = 0.1051837808310392

$data[17,0] = "MuSCs"
$data[17,1] = "Tgfb2"
$data[17,2] = "Tgfbr2"
$data[17,3] = "Inflammatory-Mac"
$data[17,4] = 3
$data[17,5] = 1
$data[17,6] = 7.520146
$data[17,7] = 22.560438
$data[17,8] = 0.324037650203293
$data[17,9] = 0.3240376502032931
$data[17,10] = 3
$data[17,11] = 1
$data[17,12] = 36.72715
$data[17,13] = 110.18145
$data[17,14] = 0.2143930638214748
$data[17,15] = 0.2143930638214748
$data[17,16] = 276.1935301639
$data[17,17] = 2485.7417714751
$data[17,18] = 0.06947142462059532
$data[17,19] = 0.06947142462059533

$data[18,0] = "MuSCs"
$data[18,1] = "Tgfb2"
$data[18,2] = "Tgfbr2"
$data[18,3] = "MuSCs"
$data[18,4] = 3
$data[18,5] = 1
$data[18,6] = 7.520146
$data[18,7] = 22.560438
$data[18,8] = 0.324037650203293
$data[18,9] = 0.3240376502032931
$data[18,10] = 3
$data[18,11] = 1
$data[18,12] = 7.381512666666667
$data[18,13] = 22.144538
$data[18,14] = 0.04308924368603855
$data[18,15] = 0.04308924368603857
$data[18,16] = 55.51005295418266
$data[18,17] = 499.590476587644
$data[18,18] = 0.01396253727306101
$data[18,19] = 0.01396253727306102

$data[19,0] = "MuSCs"
$data[19,1] = "Tgfb2"
$data[19,2] = "Tgfbr2"
$data[19,3] = "Resolving-Mac"
$data[19,4] = 3
$data[19,5] = 1
$data[19,6] = 7.520146
$data[19,7] = 22.560438
$data[19,8] = 0.324037650203293
$data[19,9] = 0.3240376502032931
$data[19,10] = 3
$data[19,11] = 1
$data[19,12] = 28.60040066666667
$data[19,13] = 85.801202
$data[19,14] = 0.1669535350673389
$data[19,15] = 0.1669535350673389
$data[19,16] = 215.0791886718307
$data[19,17] = 1935.712698046476
$data[19,18] = 0.05409923119635358
$data[19,19] = 0.05409923119635359

$ws.Range("A2:T21").Value = $data

Write-Output "done"
